$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.7484714652033114
$ws.Range("C2").Value2 = 33.8482574815652
$ws.Range("D2").Value2 = 1668.793212724136
$ws.Range("E2").Value2 = 199521.2751395922
$ws.Range("F2").Value2 = 18828212.05595672
$ws.Range("G2").Value2 = 37785982911.24001
$ws.Range("H2").Value2 = 1410165815743714
$ws.Range("B3").Value2 = 0.6998300239619539
$ws.Range("C3").Value2 = 29.27638233146019
$ws.Range("D3").Value2 = 513.5749039887279
$ws.Range("E3").Value2 = 50680.78356400538
$ws.Range("F3").Value2 = 4647295.854732187
$ws.Range("G3").Value2 = 12007879160.7819
$ws.Range("H3").Value2 = 587010965309421.9
$ws.Range("B4").Value2 = 0.699491068333178
$ws.Range("C4").Value2 = 34.41017156535496
$ws.Range("D4").Value2 = 1822.021398519945
$ws.Range("E4").Value2 = 211562.098164686
$ws.Range("F4").Value2 = 19667379.56254216
$ws.Range("G4").Value2 = 38798178024.42425
$ws.Range("H4").Value2 = 1433013179791633
$ws.Range("B5").Value2 = 0.5703222791555561
$ws.Range("C5").Value2 = 30.27576925938623
$ws.Range("D5").Value2 = 1018.038395667763
$ws.Range("E5").Value2 = 160237.4993842074
$ws.Range("F5").Value2 = 16875825.74089928
$ws.Range("G5").Value2 = 35937620393.91891
$ws.Range("H5").Value2 = 1396327342103562
$ws.Range("B6").Value2 = 0.8877020412991593
$ws.Range("C6").Value2 = 28.19110881149735
$ws.Range("D6").Value2 = 621.0127084121173
$ws.Range("E6").Value2 = 57558.21492955724
$ws.Range("F6").Value2 = 4791721.657584624
$ws.Range("G6").Value2 = 10554398124.21799
$ws.Range("H6").Value2 = 535779194320910.2
$ws.Range("B7").Value2 = 0.6646327975475445
$ws.Range("C7").Value2 = 28.91359804756351
$ws.Range("D7").Value2 = 54.94861344143851
$ws.Range("E7").Value2 = 4335.482842240131
$ws.Range("F7").Value2 = 482293.4133794098
$ws.Range("G7").Value2 = 1292964334.799403
$ws.Range("H7").Value2 = 66193307142749.42
$ws.Range("B8").Value2 = 0.8201838140663581
$ws.Range("C8").Value2 = 29.80785785109789
$ws.Range("D8").Value2 = 206.8158176423964
$ws.Range("E8").Value2 = 27775.5852990267
$ws.Range("F8").Value2 = 3089957.108081071
$ws.Range("G8").Value2 = 8139888003.539444
$ws.Range("H8").Value2 = 412188113426947.1
$ws.Range("B9").Value2 = 0.6830674722970614
$ws.Range("C9").Value2 = 29.34372483835668
$ws.Range("D9").Value2 = 596.269918181964
$ws.Range("E9").Value2 = 85253.8510114145
$ws.Range("F9").Value2 = 9156465.72864178
$ws.Range("G9").Value2 = 22094331867.46754
$ws.Range("H9").Value2 = 982541586209477.8
$ws.Range("B10").Value2 = 0.5485796179552151
$ws.Range("C10").Value2 = 28.83278660104944
$ws.Range("D10").Value2 = 27.26031922785506
$ws.Range("E10").Value2 = 27.04544196444837
$ws.Range("F10").Value2 = 28.21125293898653
$ws.Range("G10").Value2 = 1261.12908964024
$ws.Range("H10").Value2 = 84539325.02215943